$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 383.8889
$ws.Range("I5").Value = 181.875
$ws.Range("K5").Value = 181.875
$ws.Range("M5").Value = -66.875

$ws.Range("H19").Value = 3767.8
$ws.Range("I19").Value = 2961.3333
$ws.Range("K19").Value = 2961.3333
$ws.Range("M19").Value = -2786.3333

$ws.Range("H43").Value = 1604.125
$ws.Range("I43").Value = 1518.25
$ws.Range("K43").Value = 1518.25
$ws.Range("M43").Value = -1449.25

$ws.Range("H99").Value = 446.375
$ws.Range("I99").Value = 446.375
$ws.Range("K99").Value = 1339.125
$ws.Range("M99").Value = 158.875

$ws.Range("H112").Value = 1236.15
$ws.Range("J112").Value = 1298.5
$ws.Range("L112").Value = 3895.5
$ws.Range("N112").Value = -6111.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3713.1904
$ws.Range("I32").Value = 2677.4915
$ws.Range("K32").Value = 2677.4915
$ws.Range("M32").Value = -2390.4915

$ws.Range("H122").Value = 4588.16
$ws.Range("I122").Value = 4196.1816
$ws.Range("K122").Value = 12588.5448
$ws.Range("M122").Value = -10138.5448

$ws.Range("H132").Value = 2817.5
$ws.Range("I132").Value = 2955.5454
$ws.Range("J132").Value = 1299
$ws.Range("K132").Value = 8866.636200000001
$ws.Range("L132").Value = 3897
$ws.Range("M132").Value = -6336.636200000001
$ws.Range("N132").Value = -8957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2662.0698
$ws.Range("I31").Value = 2464.5151
$ws.Range("J31").Value = 3314
$ws.Range("K31").Value = 2464.5151
$ws.Range("L31").Value = 3314
$ws.Range("M31").Value = -2169.5151
$ws.Range("N31").Value = -3904

$ws.Range("H34").Value = 2662.0698
$ws.Range("I34").Value = 2464.5151
$ws.Range("J34").Value = 3314
$ws.Range("K34").Value = 2464.5151
$ws.Range("L34").Value = 3314
$ws.Range("M34").Value = -2262.5151
$ws.Range("N34").Value = -3718

$ws.Range("H70").Value = 29999
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 29999
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H105").Value = 1948.4286
$ws.Range("I105").Value = 1959.8462
$ws.Range("K105").Value = 1959.8462
$ws.Range("M105").Value = -212.8462

$ws.Range("H107").Value = 914.1842
$ws.Range("I107").Value = 542.03845
$ws.Range("J107").Value = 1720.5
$ws.Range("K107").Value = 542.03845
$ws.Range("L107").Value = 1720.5
$ws.Range("M107").Value = 1377.96155
$ws.Range("N107").Value = -5560.5

$ws.Range("H132").Value = 5825.378
$ws.Range("I132").Value = 6357.108
$ws.Range("K132").Value = 19071.324
$ws.Range("M132").Value = -16541.324

$ws.Range("H134").Value = 7504.38
$ws.Range("I134").Value = 7949.3955
$ws.Range("J134").Value = 4770.7144
$ws.Range("K134").Value = 23848.1865
$ws.Range("L134").Value = 14312.1432
$ws.Range("M134").Value = -21313.1865
$ws.Range("N134").Value = -19382.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 63634.5
$ws.Range("I4").Value = 91012.91
$ws.Range("J4").Value = 3402
$ws.Range("K4").Value = 273038.73
$ws.Range("L4").Value = 10206
$ws.Range("M4").Value = -272926.73
$ws.Range("N4").Value = -10430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 44225.2
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H70").Value = 5168.8423
$ws.Range("I70").Value = 5307.125
$ws.Range("J70").Value = 5068.273
$ws.Range("K70").Value = 5307.125
$ws.Range("L70").Value = 5068.273
$ws.Range("M70").Value = -5037.125
$ws.Range("N70").Value = -5608.273

$ws.Range("H73").Value = 5168.8423
$ws.Range("I73").Value = 5307.125
$ws.Range("J73").Value = 5068.273
$ws.Range("K73").Value = 5307.125
$ws.Range("L73").Value = 5068.273
$ws.Range("M73").Value = -4371.125
$ws.Range("N73").Value = -6940.273

$ws.Range("H81").Value = 44225.2
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 44225.2
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H132").Value = 4873.3105
$ws.Range("I132").Value = 4610.3184
$ws.Range("J132").Value = 5699.857
$ws.Range("K132").Value = 13830.9552
$ws.Range("L132").Value = 17099.571
$ws.Range("M132").Value = -11300.9552
$ws.Range("N132").Value = -22159.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2751.2222
$ws.Range("I16").Value = 2883.9048
$ws.Range("K16").Value = 2883.9048
$ws.Range("M16").Value = -2713.9048

$ws.Range("H40").Value = 11599.743
$ws.Range("I40").Value = 11649.883
$ws.Range("J40").Value = 11258.8
$ws.Range("K40").Value = 11649.883
$ws.Range("L40").Value = 11258.8
$ws.Range("M40").Value = -11513.883
$ws.Range("N40").Value = -11530.8

$ws.Range("H46").Value = 3115
$ws.Range("I46").Value = 2514.2856
$ws.Range("J46").Value = 3395.3333
$ws.Range("K46").Value = 2514.2856
$ws.Range("L46").Value = 3395.3333
$ws.Range("M46").Value = -2326.2856
$ws.Range("N46").Value = -3771.3333

$ws.Range("H55").Value = 1031.0526
$ws.Range("I55").Value = 565.3
$ws.Range("J55").Value = 1548.5555
$ws.Range("K55").Value = 565.3
$ws.Range("L55").Value = 1548.5555
$ws.Range("M55").Value = -392.3
$ws.Range("N55").Value = -1894.5555

$ws.Range("H82").Value = 611.7411499999999
$ws.Range("I82").Value = 540.35895
$ws.Range("J82").Value = 1407.1428
$ws.Range("K82").Value = 540.35895
$ws.Range("L82").Value = 1407.1428
$ws.Range("M82").Value = -179.35895
$ws.Range("N82").Value = -2129.1428

$ws.Range("H85").Value = 611.7411499999999
$ws.Range("I85").Value = 540.35895
$ws.Range("J85").Value = 1407.1428
$ws.Range("K85").Value = 540.35895
$ws.Range("L85").Value = 1407.1428
$ws.Range("M85").Value = 707.64105
$ws.Range("N85").Value = -3903.1428

$ws.Range("H136").Value = 6468.636
$ws.Range("I136").Value = 4910.7334
$ws.Range("J136").Value = 13479.2
$ws.Range("K136").Value = 14732.2002
$ws.Range("L136").Value = 40437.60000000001
$ws.Range("M136").Value = -12182.2002
$ws.Range("N136").Value = -45537.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4552583.5
$ws.Range("I81").Value = 6997298.5
$ws.Range("J81").Value = 12398.143
$ws.Range("K81").Value = 13994597
$ws.Range("L81").Value = 24796.286
$ws.Range("M81").Value = -13993536
$ws.Range("N81").Value = -26918.286

$ws.Range("H84").Value = 4552583.5
$ws.Range("I84").Value = 6997298.5
$ws.Range("J84").Value = 12398.143
$ws.Range("K84").Value = 69972985
$ws.Range("L84").Value = 123981.43
$ws.Range("M84").Value = -69967681
$ws.Range("N84").Value = -134589.43

$ws.Range("H122").Value = 2454.034
$ws.Range("I122").Value = 2842.4062
$ws.Range("J122").Value = 1993.7407
$ws.Range("K122").Value = 8527.2186
$ws.Range("L122").Value = 5981.2221
$ws.Range("M122").Value = -6077.2186
$ws.Range("N122").Value = -10881.2221

$ws.Range("H126").Value = 7489.0415
$ws.Range("I126").Value = 6813.2104
$ws.Range("J126").Value = 10057.2
$ws.Range("K126").Value = 20439.6312
$ws.Range("L126").Value = 30171.6
$ws.Range("M126").Value = -17969.6312
$ws.Range("N126").Value = -35111.60000000001

$ws.Range("H132").Value = 1434.862
$ws.Range("I132").Value = 1081.1305
$ws.Range("K132").Value = 3243.3915
$ws.Range("M132").Value = -713.3914999999997
